# Update "想去人数" (column F) counts for various events across sheets.
# Mapping of event name -> new F value (applies wherever the event appears).

$updates = @{
    "杭州·代号鸢Only——绣衣楼过大年"                                             = 1267
    "杭州·冒险家次元动漫展"                                                       = 2051
    "杭州·凡多姆海威降雪宴会·黑执事ONLY"                                          = 178
    "杭州·温馨国乙only"                                                           = 426
    "杭州·1.20新春国乙only"                                                       = 516
    "杭州·动漫迷城嘉年华"                                                         = 133
    "杭州·造梦探险家二次元同好会"                                                 = 80
    "杭州·春季任天堂同好会ONLY1.0"                                                = 163
    "杭州.第32届 中二病 原神x星穹only"                                            = 793
    "杭州·VOCALOID ONLY"                                                          = 54
    "杭州·文豪野犬舞会ONLY:横滨晚宴（取消）"                                      = 849
    "杭州·AP动漫游戏嘉年华"                                                       = 4133
    "杭州·樱之弦世界动漫游戏博览会（取消）"                                       = 2733
    "杭州·浙江蔚蓝档案only"                                                       = 848
    "杭州·偶像梦幻祭ONLY"                                                         = 622
    "杭州·第34届中二病动漫游戏展"                                                 = 720
    "杭州·6th YH樱花动漫游戏文化节"                                               = 1375
    "杭州·杭州灵能百分百only"                                                     = 288
    "杭州·异次结界动漫嘉年华"                                                     = 85
    "杭州·伊藤润二官方快闪店 限定特典礼包"                                        = 116
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value()
        if ($null -eq $name) { continue }

        if ($updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
